$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the selection/view state on the two pre-existing sheets that the
#    diff touches before we add & activate the new sheet.
# ---------------------------------------------------------------------------

# Sheet5 currently holds the "tabSelected" flag (activeTab=4). After the
# edit it is no longer the active tab and its selection moves to K2.
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Activate()
$ws5.Range("K2").Select()

# Sheet13 keeps its active cell (H6) but the selected range grows to A1:J6.
$ws13 = $wb.Worksheets.Item("Sheet13")
$ws13.Activate()
$ws13.Range("A1:J6").Select()

# ---------------------------------------------------------------------------
# 2. Add the new Sheet14 at the end of the workbook (after Sheet13). This
#    also makes it the active sheet/tab, matching activeTab="13".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws14 = $wb.Worksheets.Add($null, $lastSheet)
$ws14.Name = "Sheet14"

# Header row (row 1) - same headers as Sheet13
$ws14.Range("A1").Value = "source"
$ws14.Range("B1").Value = "author"
$ws14.Range("C1").Value = "table_name"
$ws14.Range("D1").Value = "dimensions"
$ws14.Range("E1").Value = "dimension_levels_text"
$ws14.Range("F1").Value = "dimension_levels_code"
$ws14.Range("G1").Value = "unit"
$ws14.Range("H1").Value = "interval"
$ws14.Range("I1").Value = "series_name"
$ws14.Range("J1").Value = "table_code"

# Row 2
$ws14.Range("A2").Value = "umar"
$ws14.Range("B2").Value = "mz"
$ws14.Range("C2").Value = "dfg"
$ws14.Range("D2").Value = "dff"
$ws14.Range("E2").Value = 234
$ws14.Range("F2").Value = 212
$ws14.Range("G2").Value = "%"
$ws14.Range("H2").Value = "M"
$ws14.Range("I2").Value = "tesx"
$ws14.Range("J2").Value = "MZ001"

# Row 3
$ws14.Range("A3").Value = "umar"
$ws14.Range("B3").Value = "mz"
$ws14.Range("C3").Value = "dfg"
$ws14.Range("D3").Value = "dff"
$ws14.Range("E3").Value = 1123
$ws14.Range("F3").Value = "12--de32"
$ws14.Range("G3").Value = "%"
$ws14.Range("H3").Value = "M"
$ws14.Range("I3").Value = "sdt"
$ws14.Range("J3").Value = "MZ001"

# Row 4
$ws14.Range("A4").Value = "umar"
$ws14.Range("B4").Value = "mz"
$ws14.Range("C4").Value = "dfg"
$ws14.Range("D4").Value = "dff"
$ws14.Range("E4").Value = 1123
$ws14.Range("F4").Value = "93-B"
$ws14.Range("G4").Value = "%"
$ws14.Range("H4").Value = "M"
$ws14.Range("I4").Value = "sdt"
$ws14.Range("J4").Value = "MZ001"

# Row 5
$ws14.Range("A5").Value = "umar"
$ws14.Range("B5").Value = "mz"
$ws14.Range("C5").Value = "hgf"
$ws14.Range("D5").Value = "dff"
$ws14.Range("E5").Value = 1123
$ws14.Range("F5").Value = "testiram--123"
$ws14.Range("G5").Value = "%"
$ws14.Range("H5").Value = "A"
$ws14.Range("I5").Value = "sdt"
$ws14.Range("J5").Value = "MZ002"

# Row 6
$ws14.Range("A6").Value = "umar"
$ws14.Range("B6").Value = "mz"
$ws14.Range("C6").Value = "hgf"
$ws14.Range("D6").Value = "dff"
$ws14.Range("E6").Value = 1123
$ws14.Range("F6").Value = "B,D"
$ws14.Range("G6").Value = "%"
$ws14.Range("H6").Value = "A"
$ws14.Range("I6").Value = "sdt"
$ws14.Range("J6").Value = "MZ002"

# Selection on the new sheet: active cell F6, matching the diff.
$ws14.Range("F6").Select()
